$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 145, pushing the existing rows 145-156 down to 146-157
$ws.Rows("145").Insert()

# Populate the newly inserted row 145 with the new weekly data point
$ws.Cells.Item(145, 1).Value = 10
$ws.Cells.Item(145, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(145, 3).Value = "La Araucanía"
$ws.Cells.Item(145, 4).Value = 44931
$ws.Cells.Item(145, 5).Value = 9
$ws.Cells.Item(145, 6).Value = 100114002
$ws.Cells.Item(145, 7).Value = "Camote"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 80
$ws.Cells.Item(145, 11).Value = 24000
$ws.Cells.Item(145, 12).Value = 24000
$ws.Cells.Item(145, 13).Value = 24000
$ws.Cells.Item(145, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(145, 15).Value = "Perú"
$ws.Cells.Item(145, 16).Value = 1200
$ws.Cells.Item(145, 17).Value = 20
$ws.Cells.Item(145, 18).Value = "Hortaliza"
